# ============================================================================
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" between "2021-Q4" and "总计", built as
#    a duplicate of "2021-Q4" (so header/row styling matches exactly), then
#    overwritten with the 2022-Q1 fund holdings data (21 rows).
# 2. Update the "总计" (totals) sheet: insert a new top data row for
#    2022-Q1 (21 funds, 8.48 亿元), shifting the existing 2021-Q4 / 2021-Q3
#    rows down and renumbering the index column.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Helper: force a numeric-looking (or otherwise) string into the cell as TEXT
# rather than letting Excel auto-convert it to a number. We do this by
# writing a text-returning formula into a scratch cell, copying it, and
# pasting *values only* into the destination - the pasted value keeps its
# string type without touching the destination cell's existing style.
# ----------------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, [string]$text) {
    $scratch = $ws.Cells.Item(1, 200)
    $escaped = $text -replace '"', '""'
    $scratch.Formula = '="' + $escaped + '"'
    $scratch.Copy()
    $dst = $ws.Cells.Item($row, $col)
    $dst.PasteSpecial(-4163)
    $scratch.ClearContents()
}

# ----------------------------------------------------------------------------
# Step 1: build the new "2022-Q1" sheet by duplicating "2021-Q4" (sheet #2)
# placed right before "总计" (currently sheet #3). The duplicate keeps the
# exact same header/index-column styling (s=2 bold+border+centered).
# ----------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item(2)
$wsTotalOld = $wb.Worksheets.Item(3)
$wsQ4.Copy($wsTotalOld)

$newWs = $wb.Worksheets.Item(3)
$newWs.Name = "2022-Q1"

# The source (2021-Q4) only has 12 data rows (rows 2-13); we need 21 data
# rows (rows 2-22). Extend the styled index column (A) down through row 22
# by copying the format already on A13.
$fmtSrc = $newWs.Range("A13")
$fmtSrc.Copy()
$newWs.Range("A14:A22").PasteSpecial(-4122)

# ----------------------------------------------------------------------------
# Step 2: overwrite the header row (B1:H1)
# ----------------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $col = 2 + $i
    Set-TextValue $newWs 1 $col $headers[$i]
}

# ----------------------------------------------------------------------------
# Step 3: fill in the 21 fund rows (rows 2-22)
# Columns: idx(A, number), code(B, text), name(C, text), scale(D, text),
#          position(E, text), ratio(F, text), value(G, text), rank(H, number)
# ----------------------------------------------------------------------------
$q1Data = @(
    ,@(0, "001643", "汇丰晋信智造先锋股票A", "29.09", "92.99", "4.93", "1.4341", 3)
    ,@(1, "008120", "万家自主创新混合A", "31.86", "93.78", "4.29", "1.3668", 10)
    ,@(2, "000697", "汇添富移动互联股票", "24.17", "83.18", "5.56", "1.3439", 2)
    ,@(3, "010003", "景顺长城电子信息产业股票A", "23.62", "93.46", "4.32", "1.0204", 10)
    ,@(4, "010755", "博道睿见一年持有期混合", "6.51", "93.26", "9.71", "0.6321", 1)
    ,@(5, "001644", "汇丰晋信智造先锋股票C", "10.91", "92.99", "4.93", "0.5379", 3)
    ,@(6, "506001", "万家科创板 2 年定期开放混合型证券投资基金", "12.84", "98.14", "4.02", "0.5162", 4)
    ,@(7, "010004", "景顺长城电子信息产业股票C", "7.66", "93.46", "4.32", "0.3309", 10)
    ,@(8, "008633", "万家科技创新混合A", "3.75", "93.27", "7.28", "0.2730", 9)
    ,@(9, "013123", "汇添富精选核心优势一年持有混合A", "6.15", "66.61", "4.11", "0.2528", 6)
    ,@(10, "540010", "汇丰晋信科技先锋股票", "4.37", "94.62", "5.77", "0.2521", 7)
    ,@(11, "008121", "万家自主创新混合C", "2.80", "93.78", "4.29", "0.1201", 10)
    ,@(12, "260111", "景顺长城公司治理混合", "3.08", "91.97", "3.47", "0.1069", 9)
    ,@(13, "008634", "万家科技创新混合C", "1.18", "93.27", "7.28", "0.0859", 9)
    ,@(14, "008533", "惠升惠兴混合A", "3.30", "26.14", "1.96", "0.0647", 3)
    ,@(15, "000965", "汇丰晋信新动力混合", "0.99", "91.64", "4.99", "0.0494", 3)
    ,@(16, "011077", "汇丰晋信创新先锋股票", "0.95", "94.42", "4.36", "0.0414", 10)
    ,@(17, "007152", "诺德策略精选混合", "0.42", "93.32", "5.29", "0.0222", 5)
    ,@(18, "002772", "光大保德信产业新动力灵活配置混合", "0.27", "90.63", "4.75", "0.0128", 9)
    ,@(19, "013124", "汇添富精选核心优势一年持有混合C", "0.30", "66.61", "4.11", "0.0123", 6)
    ,@(20, "008534", "惠升惠兴混合C", "0.01", "26.14", "1.96", "0.0002", 3)
)

foreach ($row in $q1Data) {
    $idxVal = [int]$row[0]
    $rankVal = [int]$row[7]
    $r = $idxVal + 2
    $newWs.Cells.Item($r, 1).Value = $idxVal
    Set-TextValue $newWs $r 2 $row[1]
    Set-TextValue $newWs $r 3 $row[2]
    Set-TextValue $newWs $r 4 $row[3]
    Set-TextValue $newWs $r 5 $row[4]
    Set-TextValue $newWs $r 6 $row[5]
    Set-TextValue $newWs $r 7 $row[6]
    $newWs.Cells.Item($r, 8).Value = $rankVal
}

# ----------------------------------------------------------------------------
# Step 4: update the "总计" (totals) sheet - insert a new row for 2022-Q1
# at the top of the data (row 2), pushing 2021-Q4 / 2021-Q3 down.
# ----------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(4)
$wsTotal.Rows.Item(2).Insert()

# The freshly inserted row inherits formatting from the row above (B1:D1,
# which is bold/bordered) - strip that so B2:D2 end up unstyled, matching
# the plain data rows elsewhere in this sheet.
$wsTotal.Range("B2:D2").ClearFormats()

# Give the new A2 the same styling (s=2) as the other index cells by
# copying the format from A3 (the row that used to be A2 before the insert).
$idxFmtSrc = $wsTotal.Range("A3")
$idxFmtSrc.Copy()
$wsTotal.Range("A2").PasteSpecial(-4122)

# Renumber the index column: new row2=0, row3=1 (was 0), row4=2 (was 1)
$wsTotal.Cells.Item(2, 1).Value = 0
$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(4, 1).Value = 2

# Fill in the new 2022-Q1 totals row
Set-TextValue $wsTotal 2 2 "2022-Q1"
$wsTotal.Cells.Item(2, 3).Value = 21
$wsTotal.Cells.Item(2, 4).Value = 8.48
